$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where Extent (column F) changes from "WV" to "WV, VA"
$extentRows = @(3,4,5,6,7,8,9,11,20,21,22,27,32,33,34,43,44,46,48,49,50,52,53,54,57,58,62,63,64,68,69,70,71,72,74,76)

foreach ($r in $extentRows) {
    $ws.Range("F$r").Value = "WV, VA"
}

# Rows where Extent Match? (column M) changes from "no" to "yes"
$matchRows = @(13,36)

foreach ($r in $matchRows) {
    $ws.Range("M$r").Value = "yes"
}
